# Added data for lean mass at sacrifice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("N1").Value = "Lean.Mass"

# Lean mass (mg) readings at sacrifice - one value per mouse row.
# Two mice (rows 2 and 15) have no recorded value.
$leanMass = [ordered]@{
    3  = 24120
    4  = 24110
    5  = 24800
    6  = 23210
    7  = 27170
    8  = 27780
    9  = 26930
    10 = 28680
    11 = 24700
    12 = 23240
    13 = 26080
    14 = 25960
    16 = 25310
    17 = 23520
    18 = 27450
    19 = 26090
    20 = 20560
    21 = 28280
    22 = 23310
    23 = 23170
    24 = 29530
    25 = 23340
    26 = 25940
    27 = 29420
    28 = 23160
    29 = 24760
}

$formatSourceCell = $null
foreach ($r in $leanMass.Keys) {
    $cell = $ws.Cells.Item($r, 14)
    $cell.Value = $leanMass[$r]

    if ($null -eq $formatSourceCell) {
        # Establish the format (16pt black Arial) on the first populated cell.
        $cell.Font.Name = "Arial"
        $cell.Font.Size = 16
        $cell.Font.Color = 0
        $formatSourceCell = $cell
    } else {
        # Re-use the exact same style for every other populated cell.
        $formatSourceCell.Copy()
        $cell.PasteSpecial(-4122)
    }

    $ws.Rows.Item($r).RowHeight = 20
}

$excel.CutCopyMode = 0
$ws.Range("N1").Select()
